$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 165
$ws.Range("I2").Value = 403
$ws.Range("J2").Value = 1765
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 454
$ws.Range("M2").Value = 31
$ws.Range("N2").Value = 293
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = 200
$ws.Range("T2").Value = 280
$ws.Range("U2").Value = 15
$ws.Range("V2").Value = 2694
$ws.Range("X2").Value = 2728
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 31
$ws.Range("AA2").Value = 19
